$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("E2").Value = 25.15000000000049
$ws.Range("G2").Value = 0.001405017506982964
$ws.Range("H2").Value = 0.003544821475448108
$ws.Range("I2").Value = "'"
$ws.Range("I2").Style = $ws.Range("J2").Style
$ws.Range("K2").Value = 3.957440810600899
$ws.Range("L2").Value = "[1.3252841496474757, 6.589597471554323]"
$ws.Range("M2").Value = 0.003322608166664631
$ws.Range("N2").Value = 0.003322608166664631
$ws.Range("O2").Value = -2.000052980663773
$ws.Range("P2").Value = "[-2.7547899544991576, -1.2453160068283875]"
$ws.Range("Q2").Value = [double]"3.270519293163687e-07"
$ws.Range("R2").Value = [double]"6.541038586327375e-07"
$ws.Range("S2").Value = 11.02264946243561
$ws.Range("T2").Value = "[9.53652673620304, 12.508772188668186]"
$ws.Range("W2").Value = 8.005705705705864
$ws.Range("X2").Value = 4.984684684684787
$ws.Range("Y2").Value = 11.02672672672694

# --- Row 3 ---
$ws.Range("E3").Value = 23.77000000000028
$ws.Range("G3").Value = [double]"2.510346732498281e-05"
$ws.Range("H3").Value = 0.0003555823058074605
$ws.Range("K3").Value = 5.151573597365312
$ws.Range("L3").Value = "[2.707611492916108, 7.595535701814516]"
$ws.Range("M3").Value = [double]"4.412684625032526e-05"
$ws.Range("N3").Value = [double]"8.825369250065052e-05"
$ws.Range("O3").Value = 1.478026573760965
$ws.Range("P3").Value = "[0.8742369946926551, 2.0818161528292745]"
$ws.Range("Q3").Value = [double]"2.356631694011568e-06"
$ws.Range("R3").Value = [double]"2.356631694011568e-06"
$ws.Range("S3").Value = 11.01759695260821
$ws.Range("T3").Value = "[9.494524333896749, 12.540669571319668]"
$ws.Range("W3").Value = 18.17845845845867
$ws.Range("X3").Value = 15.89425425425443
$ws.Range("Y3").Value = 20.4626626626629
